# Auto-generated edit script applying the cryptos price/volume update diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.329.98"
$ws.Range("D3").Value = "1.842.37"
$ws.Range("E3").Value = "  -0.30%  "
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "0.9986"
$c.Style = "Normal"
$ws.Range("E4").Value = "  -0.03%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "238.81"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -0.88%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "0.6302"
$c.Style = "Normal"
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.9996"
$c.Style = "Normal"
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.07523"
$c.Style = "Normal"
$ws.Range("E8").Value = "  -0.95%  "
$ws.Range("E9").Value = "  -1.23%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "24.43"
$c.Style = "Normal"
$ws.Range("E10").Value = "  -0.42%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.07688"
$c.Style = "Normal"
$ws.Range("E11").Value = "  -0.50%  "
$ws.Range("D12").Value = "1.826.90"
$ws.Range("E12").Value = "  -7.97%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "4.982"
$c.Style = "Normal"
$ws.Range("E13").Value = "  -0.09%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "0.6777"
$c.Style = "Normal"
$ws.Range("E14").Value = "  -1.20%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "0.00001053"
$c.Style = "Normal"
$ws.Range("E15").Value = "  +6.22%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "82.91"
$c.Style = "Normal"
$ws.Range("E16").Value = "  -0.01%  "
$ws.Range("D17").Value = "2.089.48"
$ws.Range("E17").Value = "  -7.72%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "6.114"
$c.Style = "Normal"
$ws.Range("E18").Value = "  -0.99%  "
$ws.Range("D19").Value = "29.366.01"
$ws.Range("E19").Value = "  -0.25%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "227.99"
$c.Style = "Normal"
$ws.Range("E20").Value = "  -1.50%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "12.42"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -0.82%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "0.9996"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -0.02%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "7.417"
$c.Style = "Normal"
$ws.Range("E23").Value = "  -2.59%  "
$ws.Range("E24").Value = "  -0.01%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "156.66"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +1.29%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "0.1387"
$c.Style = "Normal"
$ws.Range("E26").Value = "  -0.61%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "8.348"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -1.25%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "17.59"
$c.Style = "Normal"
$ws.Range("E28").Value = "  -0.56%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "1.452"
$c.Style = "Normal"
$ws.Range("E29").Value = "  -1.47%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "1.275"
$c.Style = "Normal"
$ws.Range("E30").Value = "  +1.07%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "0.05620"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -3.46%  "
$ws.Range("E32").Value = "  -0.63%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "4.016"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -0.16%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "1.826"
$c.Style = "Normal"
$ws.Range("E34").Value = "  -2.30%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "1.155"
$c.Style = "Normal"
$ws.Range("E35").Value = "  -0.45%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "0.7087"
$c.Style = "Normal"
$ws.Range("E36").Value = "  -1.29%  "
$ws.Range("E37").Value = "  -0.01%  "
$ws.Range("D38").Value = "1.239.91"
$ws.Range("E38").Value = "  -0.75%  "
$ws.Range("E39").Value = "  -0.18%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "2.762"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -1.28%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "0.9001"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -0.55%  "
$ws.Range("E43").Value = "  -0.02%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "101.85"
$c.Style = "Normal"
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "65.50"
$c.Style = "Normal"
$ws.Range("E45").Value = "  -2.66%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "0.00000000119"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +0.98%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "7.080"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -2.99%  "
$ws.Range("E48").Value = "  -0.65%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "8.938"
$c.Style = "Normal"
$ws.Range("E49").Value = "  -2.96%  "
$ws.Range("E50").Value = "  -1.40%  "
$ws.Range("E51").Value = "  -0.37%  "
